$d = $word.ActiveDocument
$wdReplaceOne = 1
$wdFindContinue = 1

# --- 1. Header block: "Proyecto:" line - merge "Orchid Cosmetics" (already
#     rendered as contiguous text by Find; spell-check run splits collapse
#     automatically once the surrounding runs share the same content) ---

# --- 2. Hora: "12:30 - 13:30" -> "11:30 - 12:43" ---
$d.Content.Find.Execute(
    "12:30 " + [char]8211 + " 13:30", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "11:30 " + [char]8211 + " 12:43", $wdReplaceOne) | Out-Null

# --- 3. Final "Cierre" time: "13:30 horas" -> "12:43 horas" ---
$d.Content.Find.Execute(
    "a las 13:30 horas", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "a las 12:43 horas", $wdReplaceOne) | Out-Null

# --- 4. Insert the new "Reparto de las tareas" bullet, right before the
#        "Conclusiones" heading paragraph, using the same bullet list
#        (numId 2) as the other "Desarrollo" bullets. Achieved by growing
#        a new paragraph right after the last existing bullet in that
#        list ("Mención de la posibilidad..."), which inherits its
#        numbering / spacing / run formatting automatically. ---
$modelIdx = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Mención de la posibilidad de futuras reuniones individuales con los miembros del equipo para detallar responsabilidades y recursos.") {
        $modelIdx = $i
        break
    }
    $i++
}
$model = $d.Paragraphs($modelIdx)
$modelRange = $model.Range
$modelRange.Collapse(0) | Out-Null   # wdCollapseEnd
$modelRange.InsertParagraphAfter()

$newPara = $d.Paragraphs($modelIdx + 1)
$newRange = $newPara.Range
$newRange.Collapse(1) | Out-Null     # wdCollapseStart
$newRange.InsertBefore("Reparto de las tareas")

# --- 5. Replace the "siguiente paso" sentence ---
$d.Content.Find.Execute(
    "la distribución de los paquetes de trabajo según el Plan de Dirección del Proyecto.",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "realizar un Product Backlog y un Sprint Backlog de la primera iteración.",
    $wdReplaceOne) | Out-Null

# --- 6. Remove "Acciones pendientes" section (heading + paragraph) ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Acciones pendientes") {
        $headPara = $p
        break
    }
}
$bodyPara = $headPara.Next()
$removeRange = $d.Range($headPara.Range.Start, $bodyPara.Range.End)
$removeRange.Delete()

# --- 7. Remove trailing empty paragraph at the very end of the document.
#        Delete the range spanning the paragraph mark of the last
#        non-empty paragraph through the end of the trailing empty one,
#        which merges them away. ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
if ($lastPara.Range.Text.Trim([char]13, [char]7) -eq "") {
    $secondLastPara = $d.Paragraphs($n - 1)
    $mergeRange = $d.Range($secondLastPara.Range.End - 1, $lastPara.Range.End)
    $mergeRange.Delete()
}
